# Update the "use case / integration testing" table on slide 14.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$sh = $s.Shapes.Item(2)
$tbl = $sh.Table

# Row 3 ("2" - Test case for GPA calculation)
$tbl.Cell(3, 5).Shape.TextFrame.TextRange.Text = "Shows calculated GPA"
$tbl.Cell(3, 6).Shape.TextFrame.TextRange.Text = "Pass"

# Row 6 ("5" - Test case for adding/removing courses -> administrator viewing all students)
$tbl.Cell(6, 2).Shape.TextFrame.TextRange.Text = "Test case for administrator viewing all students"
$tbl.Cell(6, 3).Shape.TextFrame.TextRange.Text = "Admin view all info"
$tbl.Cell(6, 4).Shape.TextFrame.TextRange.Text = "Student only sees itself"
$tbl.Cell(6, 5).Shape.TextFrame.TextRange.Text = "Administrator views all students"
$tbl.Cell(6, 6).Shape.TextFrame.TextRange.Text = "Pass"
